# The "id_scenario" column is redundant (every row has the same value) and is
# removed from both the worksheet and the table that wraps it, per the
# commit message: "unnecessary id_scenario columns removed from tables".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# There's exactly one table on the sheet, covering A1:F8 with columns:
# id_scenario, id_region, id_building_location, unit, 2020, 2021
$lo = $ws.ListObjects.Item(1)
$tableName = $lo.Name
$tableStyleName = $lo.TableStyle.Name

# Drop the table binding first so deleting the column doesn't leave a stale
# table definition behind, then remove column A (id_scenario) - this shifts
# id_region/id_building_location/unit/2020/2021 one column to the left.
$lo.Unlist()
$ws.Range("A1").EntireColumn.Delete()

# Re-create the table over the new A1:E8 range (now 5 columns) and restore
# its original name/style.
$newLo = $ws.ListObjects.Add(1, $ws.Range("A1:E8"), $null, 1)
$newLo.Name = $tableName
$newLo.TableStyle = $tableStyleName

$null = $ws.Range("B8").Select()
